# "site updated with new data" — refresh the pilot-data values on Sheet1
# and move the active selection from D9 to E9, matching the new data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# session 1 (row 2): bl_acc, bl_end
$ws.Range("B2").Value = 88.68
$ws.Range("D2").Value = 159

# session 3 (row 4): bl_acc, bl_end
$ws.Range("B4").Value = 87.5
$ws.Range("D4").Value = 72

# session 5 (row 6): tx_acc, tx_end
$ws.Range("C6").Value = 88.1
$ws.Range("E6").Value = 210

# session 6 (row 7): tx_acc
$ws.Range("C7").Value = 85.84

# session 7 (row 8): tx_acc, tx_end
$ws.Range("C8").Value = 92
$ws.Range("E8").Value = 275

# the author's cursor ended up one cell over after the refresh
$ws.Range("E9").Select()
